{"js": "// Update the two-digit \u00f7 one-digit practice problems in the (single) table.\n// Each problem cell's text is replaced in document order with the new\n// problem text taken from the commit's target content; formatting\n// (font/size run properties) is left untouched because we only rewrite\n// the cell's text value, not its XML.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Old text -> new text, in document (reading) order. Every occurrence is\n// unique in this document, so a strict positional walk (row-major,\n// left-to-right) unambiguously identifies which cell each mapping targets.\nconst replacements = [\n  \"82\u00f79=\", \"63\u00f75=\",\n  \"49\u00f74=\", \"63\u00f79=\",\n  \"49\u00f76=\", \"66\u00f77=\",\n  \"37\u00f74=\", \"89\u00f77=\",\n  \"99\u00f72=\", \"60\u00f73=\",\n  \"87\u00f73=\", \"32\u00f72=\",\n  \"38\u00f77=\", \"66\u00f77=\",\n  \"87\u00f79=\", \"89\u00f75=\",\n  \"69\u00f79=\", \"19\u00f76=\",\n  \"92\u00f76=\", \"82\u00f79=\",\n  \"49\u00f72=\", \"80\u00f78=\",\n  \"36\u00f79=\", \"23\u00f75=\",\n  \"86\u00f73=\", \"61\u00f73=\",\n  \"15\u00f72=\", \"94\u00f73=\",\n  \"38\u00f76=\", \"65\u00f73=\",\n  \"12\u00f72=\", \"89\u00f78=\",\n  \"69\u00f72=\", \"49\u00f72=\",\n  \"26\u00f75=\", \"51\u00f77=\",\n  \"68\u00f76=\", \"82\u00f79=\",\n  \"27\u00f73=\", \"75\u00f77=\",\n  \"72\u00f73=\", \"64\u00f78=\",\n  \"49\u00f73=\", \"56\u00f79=\",\n  \"83\u00f72=\", \"18\u00f72=\",\n  \"68\u00f79=\", \"53\u00f79=\",\n  \"64\u00f74=\", \"13\u00f77=\",\n];\n\nlet idx = 0;\nconst newValues = table.values.map((row) =>\n  row.map((cell) => {\n    if (idx < replacements.length && cell === replacements[idx]) {\n      const replacement = replacements[idx + 1];\n      idx += 2;\n      return replacement;\n    }\n    return cell;\n  })\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the two-digit \u00f7 one-digit practice problems in the (single) table.\n# Each problem cell's text is replaced with the new problem text from the\n# commit, matched positionally (row, column) and guarded by a check that\n# the existing cell text is the expected \"before\" value so we never touch\n# the wrong cell. Writing through Cell.Range.Text keeps the run's existing\n# formatting (font/size) intact, since only the text node is rewritten.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# row, column, old text, new text \u2014 document (reading) order.\n$edits = @(\n    @(1, 1, \"82\u00f79=\", \"63\u00f75=\"),\n    @(1, 2, \"49\u00f74=\", \"63\u00f79=\"),\n    @(1, 3, \"49\u00f76=\", \"66\u00f77=\"),\n    @(1, 4, \"37\u00f74=\", \"89\u00f77=\"),\n    @(1, 5, \"99\u00f72=\", \"60\u00f73=\"),\n    @(5, 1, \"87\u00f73=\", \"32\u00f72=\"),\n    @(5, 2, \"38\u00f77=\", \"66\u00f77=\"),\n    @(5, 3, \"87\u00f79=\", \"89\u00f75=\"),\n    @(5, 4, \"69\u00f79=\", \"19\u00f76=\"),\n    @(5, 5, \"92\u00f76=\", \"82\u00f79=\"),\n    @(9, 1, \"49\u00f72=\", \"80\u00f78=\"),\n    @(9, 2, \"36\u00f79=\", \"23\u00f75=\"),\n    @(9, 3, \"86\u00f73=\", \"61\u00f73=\"),\n    @(9, 4, \"15\u00f72=\", \"94\u00f73=\"),\n    @(9, 5, \"38\u00f76=\", \"65\u00f73=\"),\n    @(13, 1, \"12\u00f72=\", \"89\u00f78=\"),\n    @(13, 2, \"69\u00f72=\", \"49\u00f72=\"),\n    @(13, 3, \"26\u00f75=\", \"51\u00f77=\"),\n    @(13, 4, \"68\u00f76=\", \"82\u00f79=\"),\n    @(13, 5, \"27\u00f73=\", \"75\u00f77=\"),\n    @(17, 1, \"72\u00f73=\", \"64\u00f78=\"),\n    @(17, 2, \"49\u00f73=\", \"56\u00f79=\"),\n    @(17, 3, \"83\u00f72=\", \"18\u00f72=\"),\n    @(17, 4, \"68\u00f79=\", \"53\u00f79=\"),\n    @(17, 5, \"64\u00f74=\", \"13\u00f77=\")\n)\n\nforeach ($edit in $edits) {\n    $row = $edit[0]\n    $col = $edit[1]\n    $old = $edit[2]\n    $new = $edit[3]\n\n    $cell = $t.Cell($row, $col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($current -eq $old) {\n        $cell.Range.Text = $new\n    } else {\n        Write-Output (\"Skipped row={0} col={1}: expected '{2}' but found '{3}'\" -f $row, $col, $old, $current)\n    }\n}\n"}
